# "Remove group" operation: delete the students that belong to group 321
# from the Students sheet (previously rows 11-12), shifting the remaining
# rows (e.g. the row 13 student) up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Students")

$ws.Range("A11:C12").EntireRow.Delete()
